# CDC 29 Apr data
# - Row 97 (26-Apr-2020) total-cases figure gets revised upward.
# - A new row 98 (27-Apr-2020) is appended with that day's cumulative total.
# - Selection/scroll is left where the user would have ended up after typing
#   the new row (cell B97).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Revise the 26-Apr-2020 total.
$ws.Range("B97").Value = 962491

# Append the 27-Apr-2020 row, matching column A's existing date format.
$ws.Range("A98").Value = 43948
$ws.Range("A98").NumberFormat = "[$-409]d\-mmm\-yyyy;@"
$ws.Range("B98").Value = 981246

# Move the view/selection to where it ends up in the saved file.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 88
$ws.Range("B97").Select()
